$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Character" block (rows 3-6): rename two events, fill in loop-count (C) and 2D/3D (H) ---
$ws.Range("B3").Value = "walk"
$ws.Range("B4").Value = "click"
# B5 "pick up items" and B6 "uses items" stay the same

$ws.Range("C3:C6").Value = 1
$ws.Range("H3:H6").Value = "2d"

# --- "world" block (new rows 8-11), replacing the old "Environnement/ world" / "world change" rows ---
$ws.Range("A8").Value = "world"
$ws.Range("B9").Value = "level start"
$ws.Range("B10").Value = "level end"
$ws.Range("B11").Value = "door open"

$ws.Range("C9:C11").Value = 1
$ws.Range("H9:H11").Value = "2d"

# Clear the old "Environnement/ world" text, but keep its (centered/wrapped) formatting as a
# leftover blank-but-styled cell at its original row, and drop the custom row height.
$ws.Range("A15").ClearContents()
$ws.Range("A15").EntireRow.AutoFit()

# --- "system" block (rows 13-16), replacing the old rows 25-27 ---
$ws.Range("A13").Value = "system"
$ws.Range("B14").Value = "reboot"
$ws.Range("B15").Value = "end cycle"
$ws.Range("B16").Value = "respawn"

$ws.Range("C14:C16").Value = 1
$ws.Range("H14:H16").Value = "2d"

# --- Drop everything that used to live below the new bottom of the table ---
$ws.Range("A17:P27").ClearContents()

# --- Restore the selection to where the author left off ---
$ws.Range("I3").Select()
